$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 4866.3335
$ws.Range("I47").Value = 4866.3335
$ws.Range("K47").Value = 4866.3335
$ws.Range("M47").Value = -3894.3335
$ws.Range("H62").Value = 7417.8184
$ws.Range("I62").Value = 5866.3335
$ws.Range("J62").Value = 7999.625
$ws.Range("K62").Value = 5866.3335
$ws.Range("L62").Value = 7999.625
$ws.Range("M62").Value = -5242.3335
$ws.Range("N62").Value = -9247.625
$ws.Range("H64").Value = 3156.8572
$ws.Range("H65").Value = 7417.8184
$ws.Range("I65").Value = 5866.3335
$ws.Range("J65").Value = 7999.625
$ws.Range("K65").Value = 29331.6675
$ws.Range("L65").Value = 39998.125
$ws.Range("M65").Value = -26211.6675
$ws.Range("N65").Value = -46238.125
$ws.Range("H67").Value = 3156.8572
$ws.Range("H86").Value = 5833.1665
$ws.Range("I86").Value = 5499.5
$ws.Range("K86").Value = 5499.5
$ws.Range("M86").Value = -4376.5
$ws.Range("H89").Value = 5833.1665
$ws.Range("I89").Value = 5499.5
$ws.Range("K89").Value = 27497.5
$ws.Range("M89").Value = -21881.5
$ws.Range("H97").Value = 4824
$ws.Range("J97").Value = 3925
$ws.Range("L97").Value = 11775
$ws.Range("N97").Value = -12767
$ws.Range("H101").Value = 412.375
$ws.Range("I101").Value = 299.8
$ws.Range("J101").Value = 600
$ws.Range("K101").Value = 899.4000000000001
$ws.Range("L101").Value = 1800
$ws.Range("M101").Value = 722.5999999999999
$ws.Range("N101").Value = -5044
$ws.Range("H116").Value = 5682.4375
$ws.Range("I116").Value = 4553.727
$ws.Range("K116").Value = 4553.727
$ws.Range("M116").Value = -1111.727
$ws.Range("H137").Value = 3823.2222
$ws.Range("I137").Value = 1971.125
$ws.Range("J137").Value = 5304.9
$ws.Range("K137").Value = 5913.375
$ws.Range("L137").Value = 15914.7
$ws.Range("M137").Value = -3363.375
$ws.Range("N137").Value = -21014.7

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 503244.94
$ws.Range("I122").Value = 716778.5600000001
$ws.Range("J122").Value = 4999.8335
$ws.Range("K122").Value = 2150335.68
$ws.Range("L122").Value = 14999.5005
$ws.Range("M122").Value = -2147885.68
$ws.Range("N122").Value = -19899.5005
$ws.Range("H132").Value = 1900.9333
$ws.Range("I132").Value = 1712.2632
$ws.Range("K132").Value = 5136.7896
$ws.Range("M132").Value = -2606.7896

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3499.5
$ws.Range("I99").Value = 3576.3333
$ws.Range("J99").Value = 3345.8333
$ws.Range("K99").Value = 3576.3333
$ws.Range("L99").Value = 3345.8333
$ws.Range("M99").Value = -2078.3333
$ws.Range("N99").Value = -6341.8333
$ws.Range("H105").Value = 3195.4666
$ws.Range("I105").Value = 2374.6667
$ws.Range("K105").Value = 2374.6667
$ws.Range("M105").Value = -627.6667000000002

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1100
$ws.Range("I2").Value = 1100
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1100
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -987
$ws.Range("N2").ClearContents()
$ws.Range("H31").Value = 7171.222
$ws.Range("I31").Value = 3999.5
$ws.Range("K31").Value = 3999.5
$ws.Range("M31").Value = -3704.5
$ws.Range("H34").Value = 7171.222
$ws.Range("I34").Value = 3999.5
$ws.Range("K34").Value = 3999.5
$ws.Range("M34").Value = -3797.5
$ws.Range("H86").Value = 8987.125
$ws.Range("I86").Value = 7379.6
$ws.Range("K86").Value = 7379.6
$ws.Range("M86").Value = -6256.6
$ws.Range("H89").Value = 8987.125
$ws.Range("I89").Value = 7379.6
$ws.Range("K89").Value = 36898
$ws.Range("M89").Value = -31282
$ws.Range("H99").Value = 12582.606
$ws.Range("I99").Value = 7905.5386
$ws.Range("K99").Value = 7905.5386
$ws.Range("M99").Value = -6407.5386
$ws.Range("H107").Value = 619.9048
$ws.Range("I107").Value = 411.2143
$ws.Range("K107").Value = 411.2143
$ws.Range("M107").Value = 1508.7857
$ws.Range("H122").Value = 4527.778
$ws.Range("I122").Value = 4353
$ws.Range("J122").Value = 7499
$ws.Range("K122").Value = 13059
$ws.Range("L122").Value = 22497
$ws.Range("M122").Value = -10609
$ws.Range("N122").Value = -27397
$ws.Range("H126").Value = 12582.606
$ws.Range("I126").Value = 7905.5386
$ws.Range("K126").Value = 23716.6158
$ws.Range("M126").Value = -21246.6158
$ws.Range("H132").Value = 3764.1428
$ws.Range("I132").Value = 1451.3334
$ws.Range("K132").Value = 4354.0002
$ws.Range("M132").Value = -1824.0002
$ws.Range("H141").Value = 97999
$ws.Range("J141").Value = 97999
$ws.Range("L141").Value = 97999
$ws.Range("N141").Value = -108359

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 691.44446
$ws.Range("I5").Value = 641.25
$ws.Range("J5").Value = 731.6
$ws.Range("K5").Value = 1923.75
$ws.Range("L5").Value = 2194.8
$ws.Range("M5").Value = -1811.75
$ws.Range("N5").Value = -2418.8
$ws.Range("H107").Value = 501.9091
$ws.Range("J107").Value = 516.13794
$ws.Range("L107").Value = 1548.41382
$ws.Range("N107").Value = -5388.41382
$ws.Range("H117").Value = 633.3333
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 1072
$ws.Range("I132").Value = 1159.8
$ws.Range("J132").Value = 852.5
$ws.Range("K132").Value = 10438.2
$ws.Range("L132").Value = 7672.5
$ws.Range("M132").Value = -7908.199999999999
$ws.Range("N132").Value = -12732.5
$ws.Range("H135").Value = 691.44446
$ws.Range("I135").Value = 641.25
$ws.Range("J135").Value = 731.6
$ws.Range("K135").Value = 5771.25
$ws.Range("L135").Value = 6584.400000000001
$ws.Range("M135").Value = -3236.25
$ws.Range("N135").Value = -11654.4

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2956.8076
$ws.Range("I80").Value = 2722.5386
$ws.Range("J80").Value = 3191.077
$ws.Range("K80").Value = 2722.5386
$ws.Range("L80").Value = 3191.077
$ws.Range("M80").Value = -1724.5386
$ws.Range("N80").Value = -5187.077
$ws.Range("H83").Value = 2956.8076
$ws.Range("I83").Value = 2722.5386
$ws.Range("J83").Value = 3191.077
$ws.Range("K83").Value = 13612.693
$ws.Range("L83").Value = 15955.385
$ws.Range("M83").Value = -8620.692999999999
$ws.Range("N83").Value = -25939.385
$ws.Range("H122").Value = 86742.914
$ws.Range("I122").Value = 2768.4443
$ws.Range("J122").Value = 338666.34
$ws.Range("K122").Value = 8305.332900000001
$ws.Range("L122").Value = 1015999.02
$ws.Range("M122").Value = -5855.332900000001
$ws.Range("N122").Value = -1020899.02

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2279.9
$ws.Range("I7").Value = 1741.2858
$ws.Range("J7").Value = 3536.6667
$ws.Range("K7").Value = 1741.2858
$ws.Range("L7").Value = 3536.6667
$ws.Range("M7").Value = -1629.2858
$ws.Range("N7").Value = -3760.6667
$ws.Range("H40").Value = 1264.2667
$ws.Range("I40").Value = 1247.4286
$ws.Range("K40").Value = 1247.4286
$ws.Range("M40").Value = -1111.4286
$ws.Range("H61").Value = 1762.4546
$ws.Range("I61").Value = 1710.5
$ws.Range("J61").Value = 1901
$ws.Range("K61").Value = 1710.5
$ws.Range("L61").Value = 1901
$ws.Range("M61").Value = -1508.5
$ws.Range("N61").Value = -2305
$ws.Range("H68").Value = 2380
$ws.Range("H71").Value = 2380
$ws.Range("H113").Value = 1762.4546
$ws.Range("I113").Value = 1710.5
$ws.Range("J113").Value = 1901
$ws.Range("K113").Value = 1710.5
$ws.Range("L113").Value = 1901
$ws.Range("M113").Value = 459.5
$ws.Range("N113").Value = -6241
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550
$ws.Range("H126").Value = 2279.9
$ws.Range("I126").Value = 1741.2858
$ws.Range("J126").Value = 3536.6667
$ws.Range("K126").Value = 5223.857400000001
$ws.Range("L126").Value = 10610.0001
$ws.Range("M126").Value = -2753.857400000001
$ws.Range("N126").Value = -15550.0001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9999
$ws.Range("J41").Value = 9999
$ws.Range("L41").Value = 9999
$ws.Range("N41").Value = -10779
$ws.Range("H45").Value = 20373.555
$ws.Range("I45").Value = 25449.5
$ws.Range("J45").Value = 18923.285
$ws.Range("K45").Value = 25449.5
$ws.Range("L45").Value = 18923.285
$ws.Range("M45").Value = -24958.5
$ws.Range("N45").Value = -19905.285
$ws.Range("H122").Value = 9987.833000000001
$ws.Range("I122").Value = 8587.4
$ws.Range("K122").Value = 25762.2
$ws.Range("M122").Value = -23312.2
$ws.Range("H126").Value = 2148.1875
$ws.Range("I126").Value = 1942.8572
$ws.Range("J126").Value = 3585.5
$ws.Range("K126").Value = 5828.571599999999
$ws.Range("L126").Value = 10756.5
$ws.Range("M126").Value = -3358.571599999999
$ws.Range("N126").Value = -15696.5
$ws.Range("H132").Value = 970
$ws.Range("I132").Value = 965.1667
$ws.Range("K132").Value = 2895.5001
$ws.Range("L132").Value = 2895.5001
$ws.Range("M132").Value = -365.5001000000002
